$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.33"
$ws.Range("E2").Value = "'1.07%"
$ws.Range("D3").Value = "'37.62"
$ws.Range("E3").Value = "'0.90%"
$ws.Range("D4").Value = "'5.158"
$ws.Range("E4").Value = "'0.98%"
$ws.Range("D5").Value = "'0.07914"
$ws.Range("E5").Value = "'1.14%"
$ws.Range("D6").Value = "'4.419"
$ws.Range("E6").Value = "'0.78%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'-2.73%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.267"
$ws.Range("E8").Value = "'-0.13%"
$ws.Range("D9").Value = "'2.934"
$ws.Range("E9").Value = "'-3.20%"
$ws.Range("D10").Value = "'0.9207"
$ws.Range("E10").Value = "'-0.46%"
$ws.Range("D11").Value = "'0.1222"
$ws.Range("E11").Value = "'-8.30%"
$ws.Range("D12").Value = "'0.1919"
$ws.Range("E12").Value = "'-1.30%"
$ws.Range("D13").Value = "'0.09123"
$ws.Range("E13").Value = "'1.32%"
$ws.Range("D14").Value = "'0.03307"
$ws.Range("E14").Value = "'-3.93%"
$ws.Range("D15").Value = "'0.09615"
$ws.Range("E15").Value = "'-0.96%"
$ws.Range("D16").Value = "'0.001380"
$ws.Range("E16").Value = "'-0.78%"
$ws.Range("D17").Value = "'0.005827"
$ws.Range("E17").Value = "'-1.77%"
$ws.Range("D18").Value = "'3.526"
$ws.Range("E18").Value = "'-1.89%"
$ws.Range("E19").Value = "'0.86%"
$ws.Range("D20").Value = "'5.268"
$ws.Range("E20").Value = "'5.24%"
$ws.Range("E21").Value = "'-1.64%"
$ws.Range("E22").Value = "'4.06%"
$ws.Range("E23").Value = "'-0.15%"
$ws.Range("D24").Value = "'0.04368"
$ws.Range("E24").Value = "'0.81%"
$ws.Range("E25").Value = "'2.62%"
$ws.Range("D26").Value = "'0.004309"
$ws.Range("E26").Value = "'-4.75%"
$ws.Range("E27").Value = "'-9.70%"
$ws.Range("D39").Value = "'0.02149"
$ws.Range("E39").Value = "'-6.00%"
$ws.Range("D40").Value = "'0.05118"
$ws.Range("E40").Value = "'1.83%"
$ws.Range("D41").Value = "'0.007563"
$ws.Range("E41").Value = "'-1.09%"
$ws.Range("D42").Value = "'0.008972"
$ws.Range("E42").Value = "'-8.40%"
$ws.Range("D43").Value = "'0.1361"
$ws.Range("E43").Value = "'0.81%"
$ws.Range("D44").Value = "'0.002010"
$ws.Range("E44").Value = "'-2.50%"
$ws.Range("D45").Value = "'0.008621"
$ws.Range("E45").Value = "'2.05%"
$ws.Range("E46").Value = "'-0.97%"
$ws.Range("E47").Value = "'-0.04%"
$ws.Range("D48").Value = "'0.003325"
$ws.Range("E48").Value = "'10.21%"
$ws.Range("E49").Value = "'-7.75%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.04%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.04%"
